# Insert a new row above row 251 (shifts existing rows 251-301 down to 252-302)
# and populate it with the new record data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("251:251").Insert()

$ws.Range("A251").Value = 6
$ws.Range("B251").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C251").Value = 'Metropolitana'
$ws.Range("D251").Value = 45015
$ws.Range("D251").NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Range("E251").Value = 13
$ws.Range("F251").Value = 100112001
$ws.Range("G251").Value = 'Berenjena'
$ws.Range("H251").Value = 'Sin especificar'
$ws.Range("I251").Value = 'Primera'
$ws.Range("J251").Value = 580
$ws.Range("K251").Value = 5000
$ws.Range("L251").Value = 6000
$ws.Range("M251").Value = 5448
$ws.Range("N251").Value = '$/caja 60 unidades'
$ws.Range("O251").Value = 'Región de Arica y Parinacota'
$ws.Range("P251").Value = 91
$ws.Range("Q251").Value = 60
$ws.Range("R251").Value = 'Hortaliza'
